$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 48291.227
$ws.Range("I76").Value = 55284.367
$ws.Range("J76").Value = 4001.3333
$ws.Range("K76").Value = 55284.367
$ws.Range("L76").Value = 4001.3333
$ws.Range("M76").Value = -54969.367
$ws.Range("N76").Value = -4631.3333
$ws.Range("H79").Value = 48291.227
$ws.Range("I79").Value = 55284.367
$ws.Range("J79").Value = 4001.3333
$ws.Range("K79").Value = 55284.367
$ws.Range("L79").Value = 4001.3333
$ws.Range("M79").Value = -54192.367
$ws.Range("N79").Value = -6185.3333
$ws.Range("H138").Value = 1866.5294
$ws.Range("I138").Value = 1520.2812
$ws.Range("J138").Value = 2174.3057
$ws.Range("K138").Value = 4560.8436
$ws.Range("L138").Value = 6522.9171
$ws.Range("M138").Value = 579.1563999999998
$ws.Range("N138").Value = -16802.9171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1955
$ws.Range("I122").Value = 1955
$ws.Range("K122").Value = 5865
$ws.Range("M122").Value = -3415
$ws.Range("H132").Value = 4935.1714
$ws.Range("I132").Value = 5097.6206
$ws.Range("J132").Value = 4150
$ws.Range("K132").Value = 15292.8618
$ws.Range("L132").Value = 12450
$ws.Range("M132").Value = -12762.8618
$ws.Range("N132").Value = -17510

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13400
$ws.Range("I26").Value = 13400
$ws.Range("K26").Value = 13400
$ws.Range("M26").Value = -13108
$ws.Range("H64").Value = 371.85715
$ws.Range("I64").Value = 374.4
$ws.Range("J64").Value = 370.44446
$ws.Range("K64").Value = 374.4
$ws.Range("L64").Value = 370.44446
$ws.Range("M64").Value = -149.4
$ws.Range("N64").Value = -820.4444599999999
$ws.Range("H67").Value = 371.85715
$ws.Range("I67").Value = 374.4
$ws.Range("J67").Value = 370.44446
$ws.Range("K67").Value = 374.4
$ws.Range("L67").Value = 370.44446
$ws.Range("M67").Value = 405.6
$ws.Range("N67").Value = -1930.44446
$ws.Range("H86").Value = 3302.3914
$ws.Range("I86").Value = 2457.5386
$ws.Range("J86").Value = 4400.7
$ws.Range("K86").Value = 2457.5386
$ws.Range("L86").Value = 4400.7
$ws.Range("M86").Value = -1334.5386
$ws.Range("N86").Value = -6646.7
$ws.Range("H89").Value = 3302.3914
$ws.Range("I89").Value = 2457.5386
$ws.Range("J89").Value = 4400.7
$ws.Range("K89").Value = 12287.693
$ws.Range("L89").Value = 22003.5
$ws.Range("M89").Value = -6671.692999999999
$ws.Range("N89").Value = -33235.5
$ws.Range("H105").Value = 2864.1875
$ws.Range("I105").Value = 2786.6924
$ws.Range("K105").Value = 2786.6924
$ws.Range("M105").Value = -1039.6924
$ws.Range("H134").Value = 70695.13
$ws.Range("I134").Value = 171987.83
$ws.Range("J134").Value = 3166.6667
$ws.Range("K134").Value = 515963.49
$ws.Range("L134").Value = 9500.000100000001
$ws.Range("M134").Value = -513428.49
$ws.Range("N134").Value = -14570.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 472.65
$ws.Range("I22").Value = 411.69232
$ws.Range("J22").Value = 585.8570999999999
$ws.Range("K22").Value = 411.69232
$ws.Range("L22").Value = 585.8570999999999
$ws.Range("M22").Value = -61.69232
$ws.Range("N22").Value = -1285.8571
$ws.Range("H58").Value = 1244.3043
$ws.Range("I58").Value = 911.3889
$ws.Range("J58").Value = 2442.8
$ws.Range("K58").Value = 911.3889
$ws.Range("L58").Value = 2442.8
$ws.Range("M58").Value = -708.3889
$ws.Range("N58").Value = -2848.8
$ws.Range("H99").Value = 39505.035
$ws.Range("I99").Value = 85393.664
$ws.Range("K99").Value = 85393.664
$ws.Range("M99").Value = -83895.664
$ws.Range("H122").Value = 2214.7273
$ws.Range("I122").Value = 2882.4
$ws.Range("J122").Value = 1658.3334
$ws.Range("K122").Value = 8647.200000000001
$ws.Range("L122").Value = 4975.0002
$ws.Range("M122").Value = -6197.200000000001
$ws.Range("N122").Value = -9875.0002
$ws.Range("H126").Value = 39505.035
$ws.Range("I126").Value = 85393.664
$ws.Range("K126").Value = 256180.992
$ws.Range("M126").Value = -253710.992
$ws.Range("H132").Value = 1608.4286
$ws.Range("I132").Value = 1894.9678
$ws.Range("J132").Value = 1253.12
$ws.Range("K132").Value = 5684.903399999999
$ws.Range("L132").Value = 3759.36
$ws.Range("M132").Value = -3154.903399999999
$ws.Range("N132").Value = -8819.360000000001
$ws.Range("H134").Value = 1287.4667
$ws.Range("I134").Value = 1254.7693
$ws.Range("K134").Value = 3764.3079
$ws.Range("M134").Value = -1229.3079
$ws.Range("H136").Value = 1244.3043
$ws.Range("I136").Value = 911.3889
$ws.Range("J136").Value = 2442.8
$ws.Range("K136").Value = 2734.1667
$ws.Range("L136").Value = 7328.400000000001
$ws.Range("M136").Value = -184.1667000000002
$ws.Range("N136").Value = -12428.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 302.76666
$ws.Range("I2").Value = 327.80768
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 1966.84608
$ws.Range("L2").Value = 840
$ws.Range("M2").Value = -1853.84608
$ws.Range("N2").Value = -1066
$ws.Range("H119").Value = 15950
$ws.Range("I119").Value = 4500
$ws.Range("K119").Value = 13500
$ws.Range("M119").Value = -8662
$ws.Range("H131").Value = 2372.612
$ws.Range("I131").Value = 4785.8
$ws.Range("J131").Value = 936.1905
$ws.Range("K131").Value = 14357.4
$ws.Range("L131").Value = 2808.5715
$ws.Range("M131").Value = -9317.400000000001
$ws.Range("N131").Value = -12888.5715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 90914136
$ws.Range("I122").Value = 333343840
$ws.Range("J122").Value = 2994.75
$ws.Range("K122").Value = 1000031520
$ws.Range("L122").Value = 8984.25
$ws.Range("M122").Value = -1000029070
$ws.Range("N122").Value = -13884.25
$ws.Range("H132").Value = 2504.7666
$ws.Range("I132").Value = 2404.9092
$ws.Range("J132").Value = 2779.375
$ws.Range("K132").Value = 7214.7276
$ws.Range("L132").Value = 8338.125
$ws.Range("M132").Value = -4684.7276
$ws.Range("N132").Value = -13398.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22728852
$ws.Range("I7").Value = 1217.6
$ws.Range("J7").Value = 71430930
$ws.Range("K7").Value = 1217.6
$ws.Range("L7").Value = 71430930
$ws.Range("M7").Value = -1105.6
$ws.Range("N7").Value = -71431154
$ws.Range("H22").Value = 1221.8889
$ws.Range("I22").Value = 1923.75
$ws.Range("J22").Value = 660.4
$ws.Range("K22").Value = 1923.75
$ws.Range("L22").Value = 660.4
$ws.Range("M22").Value = -1628.75
$ws.Range("N22").Value = -1250.4
$ws.Range("H27").Value = 1221.8889
$ws.Range("I27").Value = 1923.75
$ws.Range("J27").Value = 660.4
$ws.Range("K27").Value = 1923.75
$ws.Range("L27").Value = 660.4
$ws.Range("M27").Value = -1816.75
$ws.Range("N27").Value = -874.4
$ws.Range("H40").Value = 2977
$ws.Range("I40").Value = 2233.8333
$ws.Range("J40").Value = 4463.3335
$ws.Range("K40").Value = 2233.8333
$ws.Range("L40").Value = 4463.3335
$ws.Range("M40").Value = -2097.8333
$ws.Range("N40").Value = -4735.3335
$ws.Range("H68").Value = 1333.5714
$ws.Range("I68").Value = 1046.4
$ws.Range("K68").Value = 1046.4
$ws.Range("M68").Value = -297.4000000000001
$ws.Range("H71").Value = 1333.5714
$ws.Range("I71").Value = 1046.4
$ws.Range("K71").Value = 5232
$ws.Range("M71").Value = -1488
$ws.Range("H126").Value = 22728852
$ws.Range("I126").Value = 1217.6
$ws.Range("J126").Value = 71430930
$ws.Range("K126").Value = 3652.8
$ws.Range("L126").Value = 214292790
$ws.Range("M126").Value = -1182.8
$ws.Range("N126").Value = -214297730
$ws.Range("H132").Value = 6162.0454
$ws.Range("I132").Value = 10323.348
$ws.Range("J132").Value = 1604.4286
$ws.Range("K132").Value = 30970.044
$ws.Range("L132").Value = 4813.2858
$ws.Range("M132").Value = -28440.044
$ws.Range("N132").Value = -9873.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 22000
$ws.Range("J123").Value = 22000
$ws.Range("L123").Value = 22000
$ws.Range("N123").Value = -31800

Write-Output "Applied all profit-column updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"